$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 204, shifting existing rows 204-214 down to 206-216.
$ws.Rows.Item(204).Resize(2).Insert()

# Row 204: new Flame Seedless / Especial entry
$ws.Cells.Item(204, 1).Value = 7
$ws.Cells.Item(204, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(204, 3).Value = "Ñuble"
$ws.Cells.Item(204, 4).Value = 45041
$ws.Cells.Item(204, 5).Value = 16
$ws.Cells.Item(204, 6).Value = "Fruta"
$ws.Cells.Item(204, 7).Value = 100109
$ws.Cells.Item(204, 8).Value = "Uva"
$ws.Cells.Item(204, 9).Value = 100109001
$ws.Cells.Item(204, 10).Value = "Uva"
$ws.Cells.Item(204, 11).Value = "Flame Seedless"
$ws.Cells.Item(204, 12).Value = "Especial"
$ws.Cells.Item(204, 13).Value = 60
$ws.Cells.Item(204, 14).Value = 12000
$ws.Cells.Item(204, 15).Value = 12000
$ws.Cells.Item(204, 16).Value = 12000
$ws.Cells.Item(204, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(204, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(204, 19).Value = 667
$ws.Cells.Item(204, 20).Value = 18

# Row 205: new Flame Seedless / Primera entry
$ws.Cells.Item(205, 1).Value = 7
$ws.Cells.Item(205, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(205, 3).Value = "Ñuble"
$ws.Cells.Item(205, 4).Value = 45041
$ws.Cells.Item(205, 5).Value = 16
$ws.Cells.Item(205, 6).Value = "Fruta"
$ws.Cells.Item(205, 7).Value = 100109
$ws.Cells.Item(205, 8).Value = "Uva"
$ws.Cells.Item(205, 9).Value = 100109001
$ws.Cells.Item(205, 10).Value = "Uva"
$ws.Cells.Item(205, 11).Value = "Flame Seedless"
$ws.Cells.Item(205, 12).Value = "Primera"
$ws.Cells.Item(205, 13).Value = 60
$ws.Cells.Item(205, 14).Value = 10000
$ws.Cells.Item(205, 15).Value = 10000
$ws.Cells.Item(205, 16).Value = 10000
$ws.Cells.Item(205, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(205, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(205, 19).Value = 556
$ws.Cells.Item(205, 20).Value = 18

# Apply date styling (style index 2 equivalent: copy format from a neighboring date cell)
$ws.Cells.Item(206, 4).Copy()
$ws.Range($ws.Cells.Item(204, 4), $ws.Cells.Item(205, 4)).PasteSpecial(-4122) # xlPasteFormats
$ws.Cells.Item(204, 4).Value = 45041
$ws.Cells.Item(205, 4).Value = 45041
$excel.CutCopyMode = 0
